$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.67
$ws.Range("B3").Value = 9.33
$ws.Range("E3").Value = 10.86
$ws.Range("E4").Value = 10.68
$ws.Range("C5").Value = 9.14
$ws.Range("D5").Value = 9.32
$ws.Range("F5").Value = 10.36
$ws.Range("J5").Value = 7.62
$ws.Range("E6").Value = 9.64
$ws.Range("G6").Value = 10.29
$ws.Range("F7").Value = 9.71
$ws.Range("H7").Value = 9.63
$ws.Range("G8").Value = 10.37
$ws.Range("I8").Value = 8.2
$ws.Range("J8").Value = 10.88
$ws.Range("H9").Value = 11.8
$ws.Range("E10").Value = 12.38
$ws.Range("H10").Value = 9.12
